$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "63÷3=21, 0" "49÷4=12, 1"
Replace-Text "44÷4=11, 0" "80÷5=16, 0"
Replace-Text "29÷8=3, 5" "63÷5=12, 3"
Replace-Text "95÷6=15, 5" "34÷8=4, 2"
Replace-Text "46÷7=6, 4" "70÷5=14, 0"
Replace-Text "56÷2=28, 0" "39÷4=9, 3"
Replace-Text "69÷6=11, 3" "24÷5=4, 4"
Replace-Text "64÷4=16, 0" "18÷6=3, 0"
Replace-Text "26÷5=5, 1" "42÷4=10, 2"
Replace-Text "63÷9=7, 0" "77÷3=25, 2"
Replace-Text "43÷5=8, 3" "17÷2=8, 1"
Replace-Text "99÷2=49, 1" "17÷5=3, 2"
Replace-Text "30÷9=3, 3" "85÷6=14, 1"
Replace-Text "80÷6=13, 2" "77÷7=11, 0"
Replace-Text "48÷5=9, 3" "14÷9=1, 5"
Replace-Text "20÷9=2, 2" "97÷8=12, 1"
Replace-Text "37÷3=12, 1" "67÷8=8, 3"
Replace-Text "60÷3=20, 0" "26÷3=8, 2"
Replace-Text "96÷6=16, 0" "56÷4=14, 0"
Replace-Text "92÷4=23, 0" "11÷4=2, 3"
Replace-Text "12÷6=2, 0" "55÷2=27, 1"
Replace-Text "67÷5=13, 2" "11÷6=1, 5"
Replace-Text "61÷7=8, 5" "20÷8=2, 4"
Replace-Text "46÷8=5, 6" "67÷2=33, 1"
Replace-Text "14÷5=2, 4" "79÷5=15, 4"
